$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Task Points" column (C) for every backlog item ---------
# Column C (header "Task Points" in C3) existed but rows 4-11 were blank.
# Give each task its point estimate.
$taskPoints = @{
    4  = 10
    5  = 10
    6  = 9
    7  = 7
    8  = 10
    9  = 8
    10 = 4
    11 = 6
}

foreach ($row in 4..11) {
    $ws.Range("C$row").Value = $taskPoints[$row]
}

# Match the new number cells' look (border + centered, like the rest of the
# row) to the existing formatting already used in column B of each row.
$ws.Range("B4:B11").Copy()
$ws.Range("C4:C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the saved view/selection state --------------------------------
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select()
